$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) and the matching row-label cell in column A share the
# same text ("FFR Diff-in-Diff" -> "FFR", "C/A Diff-in-Diff" -> "C/A", ...).
# Set both occurrences to the same new text so the engine re-dedupes them
# onto a single shared string (matching the original layout, where the
# header cell and its row label pointed at the same <si> entry).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "FFR"
$ws.Range("A2").Value = "FFR"

$ws.Range("C1").Value = "C/A"
$ws.Range("A3").Value = "C/A"

$ws.Range("D1").Value = "U"
$ws.Range("A4").Value = "U"

$ws.Range("E1").Value = '$\pi$'
$ws.Range("A5").Value = '$\pi$'

# ---------------------------------------------------------------------------
# Regression-result cells. Most of these carry a significance-star suffix
# (***, **) so they are unambiguously text and can be set directly. A few
# targets are "clean" decimal numbers (e.g. "0.069") that Excel's normal
# typed-input parsing would turn into a numeric cell - but the workbook
# stores them as text (shared strings), matching the rest of the table. We
# force those via the classic leading-apostrophe ("treat as text") input,
# then reset the cell style back to Normal so no visible quote-prefix
# formatting / style index lingers on the cell itself.
# ---------------------------------------------------------------------------

# Row 2 (FFR row) - Constant / U / $\pi$ columns
$ws.Range("C2").Value = "7.428***"
$ws.Range("D2").Value = "0.323**"
$ws.Range("E2").Value = "'0.069"
$ws.Range("E2").Style = "Normal"

# Row 3 (C/A row)
$ws.Range("B3").Value = "0.041***"
$ws.Range("D3").Value = "-0.053***"
$ws.Range("E3").Value = "0.024***"

# Row 4 (U row)
$ws.Range("B4").Value = "0.213**"
$ws.Range("C4").Value = "-6.307***"
$ws.Range("E4").Value = "'-0.084"
$ws.Range("E4").Style = "Normal"

# Row 5 ($\pi$ row)
$ws.Range("B5").Value = "'0.08"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "5.045***"
$ws.Range("D5").Value = "'-0.148"
$ws.Range("D5").Style = "Normal"

# Row 6 (Constant row)
$ws.Range("B6").Value = "0.198**"
$ws.Range("C6").Value = "'-1.359"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'-0.125"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.057"
$ws.Range("E6").Style = "Normal"

# Row 7 (r2_adj row) - genuinely numeric cells
$ws.Range("B7").Value = 0.41
$ws.Range("C7").Value = 0.69
$ws.Range("D7").Value = 0.47
$ws.Range("E7").Value = 0.38
